# UI UPDATE - ADD RAPPEL TO BUDGET RESUME
#
# The monthly budget table (Table 1) gets a "RAPPEL" (back-pay) entry for
# JUILLET, and SEPTEMBRE switches from all-zero to real consumption figures
# that include its own rappel. The running SOLD balance for SEPT/OCT/NOV/DEC
# drops accordingly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- JUIL row (row 9): add a 142-unit rappel, bump the TOTAL nombre ---
$t.Cell(9, 5).Range.Text = "142"
$t.Cell(9, 7).Range.Text = "4335"

# --- SEPT row (row 11): fill in real consumption + new rappel + totals ---
$t.Cell(11, 3).Range.Text = "4156"
$t.Cell(11, 4).Range.Text = "3 740 400,00"
$t.Cell(11, 5).Range.Text = "794"
$t.Cell(11, 6).Range.Text = "674 100,00"
$t.Cell(11, 7).Range.Text = "4950"
$t.Cell(11, 8).Range.Text = "4 414 500,00"
$t.Cell(11, 9).Range.Text = "19 779 400,00"

# --- OCT / NOV / DEC rows (12-14): carry the new SOLD balance forward ---
$t.Cell(12, 9).Range.Text = "19 779 400,00"
$t.Cell(13, 9).Range.Text = "19 779 400,00"
$t.Cell(14, 9).Range.Text = "19 779 400,00"
